# Add 6 new rows (50-55) for the multidisease case study: French, Spanish
# and Portuguese translations of the "tests" / "notifications" datasets,
# mirroring the existing English rows 48-49.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 50: French notifications -------------------------------------
$ws.Range("A50").Value = "multi_maladies_notifications"
$ws.Range("B50").Value = "linelist"
$ws.Range("C50").Value = "xlsx"
$ws.Range("D50").Value = 2
$ws.Range("E50").Value = 1
$ws.Range("F50").Value = "fr"
$ws.Range("G50").Value = "zzz"
$ws.Range("H50").Value = "national"
$ws.Range("I50").Value = "multidisease"
$ws.Range("J50").Value = "surveillance"
$ws.Range("K50").Value = "yes"
$ws.Range("L50").Value = 2024
$ws.Range("M50").Value = "Notifiable disease surveillance data in Feveria in 2024"
$ws.Range("N50").Value = "case_studies"
$ws.Range("O50").Value = "CC by-NC-SA 4.0"
$ws.Range("P50").Formula = '=CONCAT(SUBSTITUTE(I50," ",""),"_",J50,"_",G50,"_",L50)'
$ws.Range("Q50").Formula = '=CONCAT(P50,"_",B50,"_",D50,"_",E50,"_",J50,"_",L50)'

# --- Row 51: French tests ------------------------------------------------
$ws.Range("A51").Value = "multi_maladies_tests"
$ws.Range("B51").Value = "linelist"
$ws.Range("C51").Value = "csv"
$ws.Range("D51").Value = 1
$ws.Range("E51").Value = 1
$ws.Range("F51").Value = "fr"
$ws.Range("G51").Value = "zzz"
$ws.Range("H51").Value = "national"
$ws.Range("I51").Value = "multidisease"
$ws.Range("J51").Value = "surveillance"
$ws.Range("K51").Value = "yes"
$ws.Range("L51").Value = 2024
$ws.Range("M51").Value = "Test results for notifiable diseases reported in Feveria in 2024"
$ws.Range("N51").Value = "case_studies"
$ws.Range("O51").Value = "CC by-NC-SA 4.0"
$ws.Range("P51").Formula = '=CONCAT(SUBSTITUTE(I51," ",""),"_",J51,"_",G51,"_",L51)'
$ws.Range("Q51").Formula = '=CONCAT(P51,"_",B51,"_",D51,"_",E51,"_",J51,"_",L51)'

# --- Row 52: Spanish notifications ---------------------------------------
$ws.Range("A52").Value = "notificaciones_multienfermedad"
$ws.Range("B52").Value = "linelist"
$ws.Range("C52").Value = "xlsx"
$ws.Range("D52").Value = 2
$ws.Range("E52").Value = 1
$ws.Range("F52").Value = "es"
$ws.Range("G52").Value = "zzz"
$ws.Range("H52").Value = "national"
$ws.Range("I52").Value = "multidisease"
$ws.Range("J52").Value = "surveillance"
$ws.Range("K52").Value = "yes"
$ws.Range("L52").Value = 2024
$ws.Range("M52").Value = "Notifiable disease surveillance data in Feveria in 2024"
$ws.Range("N52").Value = "case_studies"
$ws.Range("O52").Value = "CC by-NC-SA 4.0"
$ws.Range("P52").Formula = '=CONCAT(SUBSTITUTE(I52," ",""),"_",J52,"_",G52,"_",L52)'
$ws.Range("Q52").Formula = '=CONCAT(P52,"_",B52,"_",D52,"_",E52,"_",J52,"_",L52)'

# --- Row 53: Spanish tests -------------------------------------------------
$ws.Range("A53").Value = "pruebas_multienfermedad"
$ws.Range("B53").Value = "linelist"
$ws.Range("C53").Value = "csv"
$ws.Range("D53").Value = 1
$ws.Range("E53").Value = 1
$ws.Range("F53").Value = "es"
$ws.Range("G53").Value = "zzz"
$ws.Range("H53").Value = "national"
$ws.Range("I53").Value = "multidisease"
$ws.Range("J53").Value = "surveillance"
$ws.Range("K53").Value = "yes"
$ws.Range("L53").Value = 2024
$ws.Range("M53").Value = "Test results for notifiable diseases reported in Feveria in 2024"
$ws.Range("N53").Value = "case_studies"
$ws.Range("O53").Value = "CC by-NC-SA 4.0"
$ws.Range("P53").Formula = '=CONCAT(SUBSTITUTE(I53," ",""),"_",J53,"_",G53,"_",L53)'
$ws.Range("Q53").Formula = '=CONCAT(P53,"_",B53,"_",D53,"_",E53,"_",J53,"_",L53)'

# --- Row 54: Portuguese notifications --------------------------------------
$ws.Range("A54").Value = "notificacoes_multidoencas"
$ws.Range("B54").Value = "linelist"
$ws.Range("C54").Value = "xlsx"
$ws.Range("D54").Value = 2
$ws.Range("E54").Value = 1
$ws.Range("F54").Value = "pt"
$ws.Range("G54").Value = "zzz"
$ws.Range("H54").Value = "national"
$ws.Range("I54").Value = "multidisease"
$ws.Range("J54").Value = "surveillance"
$ws.Range("K54").Value = "yes"
$ws.Range("L54").Value = 2024
$ws.Range("M54").Value = "Notifiable disease surveillance data in Feveria in 2024"
$ws.Range("N54").Value = "case_studies"
$ws.Range("O54").Value = "CC by-NC-SA 4.0"
$ws.Range("P54").Formula = '=CONCAT(SUBSTITUTE(I54," ",""),"_",J54,"_",G54,"_",L54)'
$ws.Range("Q54").Formula = '=CONCAT(P54,"_",B54,"_",D54,"_",E54,"_",J54,"_",L54)'

# --- Row 55: Portuguese tests ----------------------------------------------
$ws.Range("A55").Value = "testes_multidoencas"
$ws.Range("B55").Value = "linelist"
$ws.Range("C55").Value = "csv"
$ws.Range("D55").Value = 1
$ws.Range("E55").Value = 1
$ws.Range("F55").Value = "pt"
$ws.Range("G55").Value = "zzz"
$ws.Range("H55").Value = "national"
$ws.Range("I55").Value = "multidisease"
$ws.Range("J55").Value = "surveillance"
$ws.Range("K55").Value = "yes"
$ws.Range("L55").Value = 2024
$ws.Range("M55").Value = "Test results for notifiable diseases reported in Feveria in 2024"
$ws.Range("N55").Value = "case_studies"
$ws.Range("O55").Value = "CC by-NC-SA 4.0"
$ws.Range("P55").Formula = '=CONCAT(SUBSTITUTE(I55," ",""),"_",J55,"_",G55,"_",L55)'
$ws.Range("Q55").Formula = '=CONCAT(P55,"_",B55,"_",D55,"_",E55,"_",J55,"_",L55)'
